# CodeSystem-ValidAgeReason.xlsx update
# - URL changed from the "pythia" IG to the "cicada" IG
# - Date metadata value refreshed
# - A new "Jurisdiction" metadata row inserted after "Contact" (row 11),
#   pushing every following metadata row (Description .. Count) down by one

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$firstRow = 11
$lastRow  = 21

# --- Read the existing "Property"/"Value" pairs for rows 11..21 before
#     shifting anything (bottom block: Description .. Count) ---
$props = @()
$vals  = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $props += ,($ws.Cells.Item($r, 1).Value2)
    $vals  += ,($ws.Cells.Item($r, 2).Value2)
}

# --- Write that same block back starting one row lower, working from the
#     bottom up so a destination write never clobbers an unread source row ---
for ($i = $props.Length - 1; $i -ge 0; $i--) {
    $destRow = $firstRow + 1 + $i

    if ($destRow -gt $lastRow) {
        # Brand-new row beyond the old used range: clone formatting from the
        # last existing row first so the shifted cells keep the same style.
        $ws.Range("A" + $lastRow + ":B" + $lastRow).Copy()
        $ws.Range("A" + $destRow + ":B" + $destRow).PasteSpecial(-4122)
    }

    $ws.Cells.Item($destRow, 1).Value = $props[$i]

    if ($vals[$i] -eq $null) {
        $ws.Cells.Item($destRow, 2).Value = $null
    } elseif ($vals[$i] -eq "3") {
        # "Count" value ("3") looks numeric; force text so it round-trips as
        # a shared string rather than a number, matching the source type,
        # then re-apply the plain cell format (without touching the value
        # again) so the cell keeps the same style as its neighbours.
        $ws.Cells.Item($destRow, 2).NumberFormat = "@"
        $ws.Cells.Item($destRow, 2).Value = $vals[$i]
        $ws.Range("B" + $lastRow).Copy()
        $ws.Range("B" + $destRow).PasteSpecial(-4122)
    } else {
        $ws.Cells.Item($destRow, 2).Value = $vals[$i]
    }
}

# --- Insert the new "Jurisdiction" row in the now-vacated row 11, with an
#     explicit (non-blank, but empty) shared-string value cell for column B.
#     A leading quote forces Excel to store the value as text rather than
#     treating it as a blank cell; re-pasting the neighbour's format
#     afterwards (without touching the value again) restores the plain
#     style without disturbing the stored text type. ---
$ws.Cells.Item($firstRow, 1).Value = "Jurisdiction"
$ws.Cells.Item($firstRow, 2).Value = "'"
$ws.Range("B" + ($firstRow + 1)).Copy()
$ws.Range("B" + $firstRow).PasteSpecial(-4122)

# --- Update the URL value (row 2, column B) ---
$ws.Cells.Item(2, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/ValidAgeReason"

# --- Update the Date value (row 8, column B) ---
$ws.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"
